$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped values.
# Values are written as plain text (matching the sheet's existing inline-string
# convention of storing "329.19"-style numbers and "0.15%"-style percentages as
# text rather than numeric/percentage cell types), using a Text number format
# while assigning so Excel does not auto-convert them, then restoring the
# original "Normal" cell style so no formatting is altered.
$cellUpdates = @(
    @{ Cell = "D2"; Value = "329.18" }
    @{ Cell = "E2"; Value = "0.34%" }
    @{ Cell = "D3"; Value = "44.33" }
    @{ Cell = "E3"; Value = "0.98%" }
    @{ Cell = "D4"; Value = "5.595" }
    @{ Cell = "E4"; Value = "3.49%" }
    @{ Cell = "D5"; Value = "0.08080" }
    @{ Cell = "E5"; Value = "-0.25%" }
    @{ Cell = "D6"; Value = "2.027" }
    @{ Cell = "E6"; Value = "6.87%" }
    @{ Cell = "D7"; Value = "0.9516" }
    @{ Cell = "E7"; Value = "0.83%" }
    @{ Cell = "E8"; Value = "-7.36%" }
    @{ Cell = "D9"; Value = "0.1169" }
    @{ Cell = "E9"; Value = "-1.53%" }
    @{ Cell = "D10"; Value = "0.1852" }
    @{ Cell = "E10"; Value = "-2.35%" }
    @{ Cell = "D11"; Value = "10.25" }
    @{ Cell = "E11"; Value = "19.31%" }
    @{ Cell = "D12"; Value = "0.09721" }
    @{ Cell = "E12"; Value = "1.79%" }
    @{ Cell = "D13"; Value = "0.04523" }
    @{ Cell = "E13"; Value = "7.70%" }
    @{ Cell = "D14"; Value = "0.1067" }
    @{ Cell = "E14"; Value = "-0.12%" }
    @{ Cell = "D15"; Value = "0.001284" }
    @{ Cell = "E15"; Value = "0.80%" }
    @{ Cell = "D16"; Value = "0.04208" }
    @{ Cell = "E16"; Value = "-3.93%" }
    @{ Cell = "D17"; Value = "0.005890" }
    @{ Cell = "E17"; Value = "-0.85%" }
    @{ Cell = "D18"; Value = "3.375" }
    @{ Cell = "E18"; Value = "-5.14%" }
    @{ Cell = "D19"; Value = "4.316" }
    @{ Cell = "E19"; Value = "0.30%" }
    @{ Cell = "E21"; Value = "4.45%" }
    @{ Cell = "D22"; Value = "0.2506" }
    @{ Cell = "E22"; Value = "-3.80%" }
    @{ Cell = "D23"; Value = "0.001246" }
    @{ Cell = "E23"; Value = "0.38%" }
    @{ Cell = "D24"; Value = "0.004345" }
    @{ Cell = "E24"; Value = "0.61%" }
    @{ Cell = "D25"; Value = "0.0001190" }
    @{ Cell = "E25"; Value = "-4.03%" }
    @{ Cell = "E26"; Value = "-0.95%" }
    @{ Cell = "D38"; Value = "0.02673" }
    @{ Cell = "E38"; Value = "-1.25%" }
    @{ Cell = "D39"; Value = "0.05552" }
    @{ Cell = "E39"; Value = "0.58%" }
    @{ Cell = "D40"; Value = "0.007556" }
    @{ Cell = "E40"; Value = "-3.97%" }
    @{ Cell = "D41"; Value = "0.1408" }
    @{ Cell = "E41"; Value = "0.83%" }
    @{ Cell = "D42"; Value = "0.007955" }
    @{ Cell = "E42"; Value = "-18.51%" }
    @{ Cell = "D43"; Value = "0.002016" }
    @{ Cell = "E43"; Value = "-5.60%" }
    @{ Cell = "D44"; Value = "0.008392" }
    @{ Cell = "E44"; Value = "-12.88%" }
    @{ Cell = "D45"; Value = "0.00007167" }
    @{ Cell = "E45"; Value = "0.81%" }
    @{ Cell = "D46"; Value = "0.00000000750" }
    @{ Cell = "E46"; Value = "-0.66%" }
    @{ Cell = "D47"; Value = "0.004418" }
    @{ Cell = "E47"; Value = "27.14%" }
    @{ Cell = "D48"; Value = "0.002271" }
    @{ Cell = "E48"; Value = "-0.65%" }
    @{ Cell = "D49"; Value = "0.00002101" }
    @{ Cell = "E49"; Value = "-0.66%" }
    @{ Cell = "D50"; Value = "0.0002001" }
    @{ Cell = "E50"; Value = "-0.66%" }
)

foreach ($update in $cellUpdates) {
    $range = $ws.Range($update.Cell)
    $range.NumberFormat = "@"
    $range.Value = $update.Value
    $range.Style = "Normal"
}
